# Applies the 2025-07-30 violent-crime YTD data refresh.
# Updates the "2025" (column L) year-to-date totals -- and, where the
# source reclassified a handful of incidents, the "2023" (column J)
# totals -- across the Citywide Totals sheet, the By Neighborhood
# rollup, and every affected per-neighborhood sheet.

$wb = $excel.ActiveWorkbook

$sheetEdits = @(
    @{ Sheet = "Citywide Totals"; Cells = @(@{ Ref = "L2"; Value = 3896 }, @{ Ref = "L3"; Value = 4071 }, @{ Ref = "J4"; Value = 1085 }, @{ Ref = "L4"; Value = 1002 }, @{ Ref = "L5"; Value = 237 }, @{ Ref = "L6"; Value = 3526 }, @{ Ref = "J7"; Value = 15720 }, @{ Ref = "L7"; Value = 12732 }) },
    @{ Sheet = "Logan Square"; Cells = @(@{ Ref = "L6"; Value = 49 }, @{ Ref = "L7"; Value = 142 }) },
    @{ Sheet = "Austin"; Cells = @(@{ Ref = "L2"; Value = 244 }, @{ Ref = "L3"; Value = 277 }, @{ Ref = "L6"; Value = 226 }, @{ Ref = "L7"; Value = 831 }) },
    @{ Sheet = "Garfield Park"; Cells = @(@{ Ref = "L2"; Value = 165 }, @{ Ref = "L6"; Value = 190 }, @{ Ref = "L7"; Value = 596 }) },
    @{ Sheet = "Grand Crossing"; Cells = @(@{ Ref = "L3"; Value = 155 }, @{ Ref = "L4"; Value = 28 }, @{ Ref = "L6"; Value = 128 }, @{ Ref = "L7"; Value = 467 }) },
    @{ Sheet = "New City"; Cells = @(@{ Ref = "L3"; Value = 70 }, @{ Ref = "L7"; Value = 238 }) },
    @{ Sheet = "By Neighborhood"; Cells = @(@{ Ref = "L6"; Value = 104 }, @{ Ref = "L8"; Value = 831 }, @{ Ref = "L10"; Value = 80 }, @{ Ref = "L11"; Value = 210 }, @{ Ref = "L15"; Value = 93 }, @{ Ref = "L19"; Value = 368 }, @{ Ref = "L20"; Value = 317 }, @{ Ref = "L23"; Value = 135 }, @{ Ref = "L25"; Value = 70 }, @{ Ref = "L29"; Value = 703 }, @{ Ref = "L33"; Value = 596 }, @{ Ref = "L36"; Value = 164 }, @{ Ref = "L37"; Value = 467 }, @{ Ref = "L42"; Value = 408 }, @{ Ref = "L43"; Value = 97 }, @{ Ref = "L44"; Value = 91 }, @{ Ref = "L47"; Value = 93 }, @{ Ref = "L51"; Value = 157 }, @{ Ref = "L52"; Value = 256 }, @{ Ref = "L53"; Value = 142 }, @{ Ref = "L54"; Value = 264 }, @{ Ref = "L55"; Value = 119 }, @{ Ref = "J63"; Value = 104 }, @{ Ref = "L63"; Value = 42 }, @{ Ref = "L64"; Value = 85 }, @{ Ref = "L65"; Value = 238 }, @{ Ref = "L67"; Value = 451 }, @{ Ref = "L71"; Value = 36 }, @{ Ref = "L76"; Value = 193 }, @{ Ref = "L78"; Value = 160 }, @{ Ref = "L79"; Value = 333 }, @{ Ref = "L84"; Value = 123 }, @{ Ref = "L85"; Value = 669 }, @{ Ref = "L88"; Value = 141 }, @{ Ref = "L89"; Value = 178 }, @{ Ref = "L93"; Value = 69 }, @{ Ref = "L96"; Value = 134 }, @{ Ref = "J101"; Value = 15720 }, @{ Ref = "L101"; Value = 12732 }) },
    @{ Sheet = "North Lawndale"; Cells = @(@{ Ref = "L6"; Value = 105 }, @{ Ref = "L7"; Value = 451 }) },
    @{ Sheet = "South Deering"; Cells = @(@{ Ref = "L3"; Value = 46 }, @{ Ref = "L4"; Value = 5 }, @{ Ref = "L6"; Value = 29 }, @{ Ref = "L7"; Value = 123 }) },
    @{ Sheet = "Loop"; Cells = @(@{ Ref = "L6"; Value = 125 }, @{ Ref = "L7"; Value = 264 }) },
    @{ Sheet = "Englewood"; Cells = @(@{ Ref = "L2"; Value = 219 }, @{ Ref = "L3"; Value = 263 }, @{ Ref = "L7"; Value = 703 }) },
    @{ Sheet = "Chatham"; Cells = @(@{ Ref = "L2"; Value = 130 }, @{ Ref = "L3"; Value = 112 }, @{ Ref = "L6"; Value = 107 }, @{ Ref = "L7"; Value = 368 }) },
    @{ Sheet = "Irving Park"; Cells = @(@{ Ref = "L2"; Value = 38 }, @{ Ref = "L7"; Value = 91 }) },
    @{ Sheet = "River North"; Cells = @(@{ Ref = "L2"; Value = 39 }, @{ Ref = "L6"; Value = 88 }, @{ Ref = "L7"; Value = 193 }) },
    @{ Sheet = "Ashburn"; Cells = @(@{ Ref = "L2"; Value = 46 }, @{ Ref = "L7"; Value = 104 }) },
    @{ Sheet = "Humboldt Park"; Cells = @(@{ Ref = "L3"; Value = 132 }, @{ Ref = "L7"; Value = 408 }) },
    @{ Sheet = "Avondale"; Cells = @(@{ Ref = "L6"; Value = 22 }, @{ Ref = "L7"; Value = 80 }) },
    @{ Sheet = "Rogers Park"; Cells = @(@{ Ref = "L4"; Value = 17 }, @{ Ref = "L6"; Value = 48 }, @{ Ref = "L7"; Value = 160 }) },
    @{ Sheet = "Lower West Side"; Cells = @(@{ Ref = "L2"; Value = 39 }, @{ Ref = "L6"; Value = 33 }, @{ Ref = "L7"; Value = 119 }) },
    @{ Sheet = "Douglas"; Cells = @(@{ Ref = "L2"; Value = 34 }, @{ Ref = "L6"; Value = 36 }, @{ Ref = "L7"; Value = 135 }) },
    @{ Sheet = "West Ridge"; Cells = @(@{ Ref = "L3"; Value = 37 }, @{ Ref = "L7"; Value = 134 }) },
    @{ Sheet = "Roseland"; Cells = @(@{ Ref = "L3"; Value = 121 }, @{ Ref = "L6"; Value = 69 }, @{ Ref = "L7"; Value = 333 }) },
    @{ Sheet = "Near South Side"; Cells = @(@{ Ref = "L4"; Value = 10 }, @{ Ref = "L7"; Value = 85 }) },
    @{ Sheet = "Chicago Lawn"; Cells = @(@{ Ref = "L6"; Value = 82 }, @{ Ref = "L7"; Value = 317 }) },
    @{ Sheet = "Grand Boulevard"; Cells = @(@{ Ref = "L4"; Value = 11 }, @{ Ref = "L7"; Value = 164 }) },
    @{ Sheet = "West Lawn"; Cells = @(@{ Ref = "L2"; Value = 24 }, @{ Ref = "L7"; Value = 69 }) },
    @{ Sheet = "East Side"; Cells = @(@{ Ref = "L2"; Value = 24 }, @{ Ref = "L3"; Value = 35 }, @{ Ref = "L7"; Value = 70 }) },
    @{ Sheet = "Kenwood"; Cells = @(@{ Ref = "L2"; Value = 33 }, @{ Ref = "L7"; Value = 93 }) },
    @{ Sheet = "Brighton Park"; Cells = @(@{ Ref = "L4"; Value = 9 }, @{ Ref = "L7"; Value = 93 }) },
    @{ Sheet = "Belmont Cragin"; Cells = @(@{ Ref = "L4"; Value = 16 }, @{ Ref = "L7"; Value = 210 }) },
    @{ Sheet = "United Center"; Cells = @(@{ Ref = "L5"; Value = 3 }, @{ Ref = "L7"; Value = 141 }) },
    @{ Sheet = "Uptown"; Cells = @(@{ Ref = "L6"; Value = 46 }, @{ Ref = "L7"; Value = 178 }) },
    @{ Sheet = "Little Italy, UIC"; Cells = @(@{ Ref = "L2"; Value = 48 }, @{ Ref = "L7"; Value = 157 }) },
    @{ Sheet = "Hyde Park"; Cells = @(@{ Ref = "L6"; Value = 33 }, @{ Ref = "L7"; Value = 97 }) },
    @{ Sheet = "South Shore"; Cells = @(@{ Ref = "L2"; Value = 201 }, @{ Ref = "L4"; Value = 45 }, @{ Ref = "L7"; Value = 669 }) },
    @{ Sheet = "Oakland"; Cells = @(@{ Ref = "L2"; Value = 16 }, @{ Ref = "L3"; Value = 12 }, @{ Ref = "L7"; Value = 36 }) },
    @{ Sheet = "Little Village"; Cells = @(@{ Ref = "L2"; Value = 88 }, @{ Ref = "L3"; Value = 78 }, @{ Ref = "L7"; Value = 256 }) }
)

foreach ($entry in $sheetEdits) {
    $ws = $wb.Worksheets.Item($entry.Sheet)
    foreach ($cell in $entry.Cells) {
        $ws.Range($cell.Ref).Value = $cell.Value
    }
}
